$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of test/task data (row 14)
$ws.Range("A14").Value = "replace errFlag with an 8 bit (?) binary code where each bit indicates a different error?"
$ws.Range("B14").Value = "N"
$ws.Range("C14").Value = "N"
$ws.Range("D14").Value = "All"

# Move active selection to D15, matching the resulting workbook view state
$ws.Range("D15").Select()
